# PIM-6333: Fix fixtures files
#
# The "family_variant" fixture workbook had a couple of stale/incorrect
# values in its header-configuration rows (leftover "image_1" attribute
# code instead of "image", and "EAN" instead of lower-cased "ean"), plus
# a workbook that was saved with a zoomed-out tab ratio and a stale
# selection. This brings the fixture back in line with the real export
# format and tidies up the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "variant-attributes" sample values (row 2 / H2 & I2) ---
# image_1 -> image
$ws.Range("H2").Value = "color,name,image,variation_image,composition"
# EAN -> ean
$ws.Range("I2").Value = "size,ean,sku,weight"

# --- Fix the "variant-attributes_1" sample value on row 4 (H4) ---
$ws.Range("H4").Value = "name,image,variation_image,composition"

# --- Update the saved selection / active cell to E20 ---
[void]$ws.Range("E20").Select()

# --- Widen the saved tab ratio (bookViews/workbookView@tabRatio, 993/1000) ---
$wb.Windows.Item(1).TabRatio = 0.993
